$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated financial figures per row (data corrections / re-scaled units).
# A value of $null means the cell is cleared entirely (matches the source edit,
# which removes those <c> elements rather than leaving a blank value).

$row2 = [ordered]@{
    "D" = 39283
    "E" = 1517
    "F" = 1517
    "G" = 1571
    "H" = 1210
    "I" = 1210
    "J" = 0
    "K" = 246811
    "L" = 231830
    "M" = 14981
    "N" = 14978
    "O" = 3
    "P" = 6315
    "Q" = -996
    "R" = 1072
    "S" = 762
    "T" = 144
    "U" = $null
    "V" = 1728
    "W" = 3.86
    "X" = 3.08
    "Y" = 9.81
    "Z" = 0.52
    "AA" = 1547.51
    "AB" = 137.21
    "AC" = 907
    "AD" = $null
    "AE" = 11232
    "AF" = 0
    "AG" = 100
    "AH" = $null
    "AI" = 25.22
    "AJ" = 105180000
}

$row3 = [ordered]@{
    "D" = 42720
    "E" = 1228
    "F" = 1228
    "G" = 1255
    "H" = 1224
    "I" = 1223
    "J" = 0
    "K" = 272203
    "L" = 253630
    "M" = 18573
    "N" = 18570
    "O" = 3
    "P" = 8315
    "Q" = 19970
    "R" = -19923
    "S" = 459
    "T" = 2695
    "U" = $null
    "V" = 1606
    "W" = 2.88
    "X" = 2.86
    "Y" = 7.29
    "Z" = 0.47
    "AA" = 1365.57
    "AB" = 124.59
    "AC" = 797
    "AD" = 6.33
    "AE" = 10819
    "AF" = 0.47
    "AG" = 55
    "AH" = 1.09
    "AI" = 22.8
    "AJ" = 145180000
}

$row4 = [ordered]@{
    "D" = 41702
    "E" = 1168
    "F" = 1168
    "G" = 1224
    "H" = 780
    "I" = 780
    "J" = $null
    "K" = 286441
    "L" = 267639
    "M" = 18802
    "N" = 18802
    "O" = $null
    "P" = 8315
    "Q" = 13503
    "R" = -14664
    "S" = 712
    "T" = 198
    "U" = $null
    "V" = 1360
    "W" = 2.8
    "X" = 1.87
    "Y" = 4.17
    "Z" = 0.28
    "AA" = 1423.47
    "AB" = 127.34
    "AC" = 459
    "AD" = 11.21
    "AE" = 11422
    "AF" = 0.45
    "AG" = 55
    "AH" = 1.07
    "AI" = 29.36
    "AJ" = 145180000
}

$row5 = [ordered]@{
    "D" = 43220
    "E" = 640
    "F" = 640
    "G" = 2333
    "H" = 2211
    "I" = 2211
    "J" = $null
    "K" = 355208
    "L" = 334692
    "M" = 20516
    "N" = 20516
    "O" = $null
    "P" = 8315
    "Q" = 9896
    "R" = -9252
    "S" = -1301
    "T" = 378
    "U" = $null
    "V" = 1290
    "W" = 1.48
    "X" = 5.12
    "Y" = 11.25
    "Z" = 0.6899999999999999
    "AA" = 1631.37
    "AB" = 147.95
    "AC" = 1330
    "AD" = 3.95
    "AE" = 12464
    "AF" = 0.42
    "AG" = 170
    "AH" = 3.24
    "AI" = 17.82
    "AJ" = 145180000
}

$row6 = [ordered]@{
    "D" = 37926
    "E" = 1424
    "F" = 1424
    "G" = 1395
    "H" = 1066
    "I" = 1066
    "K" = 346405
    "L" = 325755
    "M" = 20650
    "N" = 20650
    "P" = 9907
    "Q" = -3326
    "R" = 2210
    "S" = 1186
    "T" = 593
    "U" = $null
    "V" = 1993
    "W" = 3.75
    "X" = 2.81
    "Y" = 5.18
    "Z" = 0.3
    "AA" = 1577.53
    "AB" = 112.33
    "AC" = 553
    "AD" = 8.34
    "AE" = 12938
    "AF" = 0.36
    "AG" = 170
    "AH" = 3.69
    "AI" = 36.15
    "AJ" = 177016189
}

$row7 = [ordered]@{
    "D" = $null
    "E" = $null
    "G" = $null
    "H" = $null
    "I" = $null
    "K" = $null
    "L" = $null
    "M" = $null
    "N" = $null
    "P" = $null
    "Q" = $null
    "R" = $null
    "S" = $null
    "T" = $null
    "U" = $null
    "W" = $null
    "X" = $null
    "Y" = $null
    "Z" = $null
    "AA" = $null
    "AC" = $null
    "AD" = $null
    "AE" = $null
    "AF" = $null
    "AG" = 191
    "AH" = 4.54
    "AI" = $null
}

$row8 = [ordered]@{
    "D" = $null
    "E" = $null
    "G" = $null
    "H" = $null
    "I" = $null
    "K" = $null
    "L" = $null
    "M" = $null
    "N" = $null
    "P" = $null
    "Q" = $null
    "R" = $null
    "S" = $null
    "T" = $null
    "U" = $null
    "W" = $null
    "X" = $null
    "Y" = $null
    "Z" = $null
    "AA" = $null
    "AC" = $null
    "AD" = $null
    "AE" = $null
    "AF" = $null
    "AG" = 197
    "AH" = 4.68
    "AI" = $null
}

$row9 = [ordered]@{
    "D" = $null
    "E" = $null
    "G" = $null
    "H" = $null
    "I" = $null
    "K" = $null
    "L" = $null
    "M" = $null
    "N" = $null
    "P" = $null
    "Q" = $null
    "R" = $null
    "S" = $null
    "T" = $null
    "U" = $null
    "W" = $null
    "X" = $null
    "Y" = $null
    "Z" = $null
    "AA" = $null
    "AC" = $null
    "AD" = $null
    "AE" = $null
    "AF" = $null
    "AG" = 215
    "AH" = 5.11
    "AI" = $null
}

$rowsData = [ordered]@{
    2 = $row2
    3 = $row3
    4 = $row4
    5 = $row5
    6 = $row6
    7 = $row7
    8 = $row8
    9 = $row9
}

foreach ($rowNum in $rowsData.Keys) {
    $cols = $rowsData[$rowNum]
    foreach ($colLetter in $cols.Keys) {
        $value = $cols[$colLetter]
        $cellRef = "$colLetter$rowNum"
        if ($null -eq $value) {
            $ws.Range($cellRef).Value = ""
        } else {
            $ws.Range($cellRef).Value = $value
        }
    }
}
